$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new data rows (Account/Customer numbers for TC 118518), appended right
# below the existing five rows. The ID-like values in columns A:D are stored
# as text in the workbook (same as every existing row), not as numbers, so
# mark the range as Text first - otherwise Excel would auto-coerce the
# all-digit strings to numeric cells. ClearFormats() afterwards drops the
# temporary Text number-format again so the new cells end up with the same
# (default/general) cell style as the rest of the sheet.
$newRows = $ws.Range("A6:D7")
$newRows.NumberFormat = "@"

$ws.Range("A6").Value = "118518"
$ws.Range("B6").Value = "1008784258"
$ws.Range("C6").Value = "17866752"
$ws.Range("D6").Value = "6020"

$ws.Range("A7").Value = "118518"
$ws.Range("B7").Value = "1008784270"
$ws.Range("C7").Value = "17866761"
$ws.Range("D7").Value = "6020"

$newRows.ClearFormats()
